$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Simple text field updates (Créditos-trabalho / Carga horária / Ativação)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Créditos-trabalho: 4", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Créditos-trabalho: 0", 2)
$d.Content.Find.Execute("Carga horária: 150 h", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Carga horária: 30 h", 2)
$d.Content.Find.Execute("Ativação: 01/01/2016", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Ativação: 01/01/2025", 2)

# ---------------------------------------------------------------------
# 2) Collapse the manual line breaks (<w:br/>) that separate the <w:t>
#    runs of text inside single runs into one contiguous <w:t>. "^l" is
#    Word's Find/Replace code for a manual line break.
# ---------------------------------------------------------------------

# English "Objetivos" paragraph (index 7): 2 segments -> 1
$pObjEn = $d.Paragraphs(7).Range
$pObjEn.Find.Execute("degree. ^l2 - Integration", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "degree. 2 - Integration", 2)

# Portuguese "Programa" paragraph (index 14): 4 segments -> 1
$pProgPt = $d.Paragraphs(14).Range
$pProgPt.Find.Execute("(P&ID).^l2  Estrutura", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "(P&ID).2  Estrutura", 2)
$pProgPt = $d.Paragraphs(14).Range
$pProgPt.Find.Execute("reciclo; ^l3  Análise", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "reciclo; 3  Análise", 2)
$pProgPt = $d.Paragraphs(14).Range
$pProgPt.Find.Execute("processos.^l4  Estudo", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "processos.4  Estudo", 2)

# English "Programa" paragraph (index 15): 4 segments -> 1
$pProgEn = $d.Paragraphs(15).Range
$pProgEn.Find.Execute("(P&ID).^l2 - The Structure", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "(P&ID).2 - The Structure", 2)
$pProgEn = $d.Paragraphs(15).Range
$pProgEn.Find.Execute("Recycle Structure of the Process^l3 - Analysis", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "Recycle Structure of the Process3 - Analysis", 2)
$pProgEn = $d.Paragraphs(15).Range
$pProgEn.Find.Execute("performance.^l4 - Industrial", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "performance.4 - Industrial", 2)
